# "Various edits based on Ed's comments."
#
# Slide 2 (index 2 in the deck) - "Metadata" -> "Memory Metadata"
# Slide 1 (index 1 in the deck) - "Main " + "Core" runs merged into one run "Main Core"
# Slide 1 - both "RF Metadata" boxes lose their (now redundant) trailing endParaRPr
# Slide 2 - "Monitoring Core" -> "Monitor"

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

# 1) "Metadata" -> "Memory Metadata" (box above "Monitoring Core" on slide 2)
$metadataShape = $s2.Shapes.Item("Rectangle 76")
$metadataShape.TextFrame.TextRange.Text = "Memory Metadata"

# 2) "Main " + "Core" (two runs, one paragraph) -> single run "Main Core" (slide 1)
$mainCoreShape = $s1.Shapes.Item("Rectangle 13")
$mainCoreTr = $mainCoreShape.TextFrame.TextRange
$mainCoreFirstPara = $mainCoreTr.Paragraphs(1, 1)
$mainCoreFirstPara.Characters(1, $mainCoreFirstPara.Length).Text = "Main Core"

# 3) & 4) "RF Metadata" boxes (slide 1) - re-typing the identical text drops the
#         stray trailing endParaRPr that no longer carries any distinct formatting.
$rfShape1 = $s1.Shapes.Item("Rectangle 15")
$rfTr1 = $rfShape1.TextFrame.TextRange
$rfTr1.Delete()
$rfTr1.Text = "RF Metadata"

$rfShape2 = $s1.Shapes.Item("Rectangle 18")
$rfTr2 = $rfShape2.TextFrame.TextRange
$rfTr2.Delete()
$rfTr2.Text = "RF Metadata"

# 5) "Monitoring Core" -> "Monitor" (slide 2)
$monitorShape = $s2.Shapes.Item("Rectangle 12")
$monitorShape.TextFrame.TextRange.Text = "Monitor"
